$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was captured for "Cilantro" at the
# Macroferia Regional de Talca. The new record is inserted as row 55
# (the data appears to be ordered with this new entry first), which
# pushes every existing row from 55 down through 129 down by one row
# (they become rows 56-130). The worksheet's used range grows from
# A1:R129 to A1:R130 as a result.
$ws.Rows("55:55").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A55").Value = 5
$ws.Range("B55").Value = 'Macroferia Regional de Talca'
$ws.Range("C55").Value = 'Maule'
$ws.Range("D55").Value = 45195
$ws.Range("E55").Value = 7
$ws.Range("F55").Value = 100112040
$ws.Range("G55").Value = 'Cilantro'
$ws.Range("H55").Value = 'Sin especificar'
$ws.Range("I55").Value = 'Primera'
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = 8000
$ws.Range("N55").Value = '$/caja 36 atados'
$ws.Range("O55").Value = 'Región del Maule'
$ws.Range("P55").Value = 222
$ws.Range("Q55").Value = 36
$ws.Range("R55").Value = 'Hortaliza'
